# B6-PowerPoint.pptx edit script
#
# 1. Three tables (on slides 14, 15 and 16) switch from the custom
#    "Table_0" style ({B18E7987-DCD1-4129-AD9D-4C8D9EE1CDE6}) to the
#    built-in table style {B19D9DF9-C451-4306-A08B-196D57DFA32A}.
# 2. The deck's theme colour scheme (currently the "Integral" / "Red Violet"
#    palette) is swapped for the plain "Office Theme" palette.

$p = $ppt.ActivePresentation

# --- 1. Re-style the three tables ------------------------------------------
$targetStyleId = "{B19D9DF9-C451-4306-A08B-196D57DFA32A}"
foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.HasTable) {
            $shp.Table.ApplyStyle($targetStyleId)
        }
    }
}

# --- 2. Swap the theme colour scheme ---------------------------------------
function RgbLong($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Office Theme colour scheme (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink)
$officeColors = @(
    "000000", "FFFFFF", "44546A", "E7E6E6",
    "5B9BD5", "ED7D31", "A5A5A5", "FFC000", "4472C4", "70AD47",
    "0563C1", "954F72"
)

$master = $p.SlideMaster
$theme = $master.Theme
$colorScheme = $theme.ThemeColorScheme
for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Item($i).RGB = RgbLong($officeColors[$i - 1])
}
